# Apply the "AddressBook" -> "GradTrak" rename described by the commit
# ("updated DG and diagrams") to the class-diagram slide.
#
# The diff only changes the text of a single run (the class-name label
# "AddressBook" inside a two-paragraph shape "AddressBook" / "Parser");
# everything else about the run (rPr, color, etc.) is left untouched.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)

        if (-not $shp.HasTextFrame) { continue }
        if (-not $shp.TextFrame.HasText) { continue }

        $tr = $shp.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)
            # Paragraph text comes back with a trailing CR (0x0D) that isn't
            # part of the actual paragraph content, so strip it before
            # comparing.
            $paraText = $para.Text.TrimEnd("`r")

            if ($paraText -eq "AddressBook") {
                # Replace just this paragraph's characters so the existing
                # run formatting (rPr/solidFill/etc.) is preserved and we
                # end up with a single run, matching the source edit.
                $len = $paraText.Length
                $chars = $para.Characters(1, $len)
                $chars.Text = "GradTrak"
            }
        }
    }
}
